$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140-194 down to 141-195.
$ws.Rows.Item(140).Insert()

# Populate the newly-inserted row 140 with the new weekly price-report record.
$ws.Cells.Item(140, 1).Value = 11
$ws.Cells.Item(140, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(140, 3).Value = "Bíobío"
$ws.Cells.Item(140, 4).Value = 44784
$ws.Cells.Item(140, 5).Value = 8
$ws.Cells.Item(140, 6).Value = "Fruta"
$ws.Cells.Item(140, 7).Value = 100108
$ws.Cells.Item(140, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(140, 9).Value = 100108005
$ws.Cells.Item(140, 10).Value = "Piña"
$ws.Cells.Item(140, 11).Value = "Caramelo"
$ws.Cells.Item(140, 12).Value = "Segunda"
$ws.Cells.Item(140, 13).Value = 220
$ws.Cells.Item(140, 14).Value = 18000
$ws.Cells.Item(140, 15).Value = 19000
$ws.Cells.Item(140, 16).Value = 18545
$ws.Cells.Item(140, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(140, 18).Value = "Ecuador"
$ws.Cells.Item(140, 19).Value = 1325
$ws.Cells.Item(140, 20).Value = 14
